# Bang trach nhiem - apply commit "Cap nhat bang trach nhiem"
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Remove the _GoBack bookmark from the "Stt" header cell (row 1, col 1) ---
$sttCell = $t.Cell(1, 1)
$sttCell.Range.Delete()
$sttCell.Range.InsertAfter("Stt")

# --- 2. Row Stt=2 (table row 3): Ghi chu cell, "thieu" -> "sai" (unique word, safe Find/Replace) ---
$d.Content.Find.Execute("thiếu", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "sai", 2) | Out-Null

# --- 3. Row Stt=3 (table row 4): rewrite Nghiep vu / Phan mem / Ghi chu cells ---
$row4col2 = $t.Cell(4, 2)
$row4col2.Range.Delete()
$row4col2.Range.InsertAfter("Tạo danh sách sản phẩm ")

$row4col4 = $t.Cell(4, 4)
$row4col4.Range.Delete()
$row4col4.Range.InsertAfter("Hiển thị thông tin sản phẩm,chi tiết sản phẩm")

$row4col5 = $t.Cell(4, 5)
$row4col5.Range.Delete()
$row4col5.Range.InsertAfter("Cập nhật trạng thái hiển thị")

# --- 4. Row Stt=4 (table row 5): Phan mem cell ---
$row5col4 = $t.Cell(5, 4)
$row5col4.Range.Delete()
$row5col4.Range.InsertAfter("Kiểm tra thông tin,thông báo,tạo hóa đơn")

# --- 5. Row Stt=5 (table row 6): Phan mem cell ---
$row6col4 = $t.Cell(6, 4)
$row6col4.Range.Delete()
$row6col4.Range.InsertAfter("Kiểm tra thông tin,thông báo,tạo")

# --- 6. Row Stt=6 (table row 7): Nguoi dung + Phan mem cells ---
$row7col3 = $t.Cell(7, 3)
$row7col3.Range.Delete()
$row7col3.Range.InsertAfter("Cung cấp thông tin chỉnh sửa ")

$row7col4 = $t.Cell(7, 4)
$row7col4.Range.Delete()
$row7col4.Range.InsertAfter("Kiểm tra thông tin và thực hiện xóa , sửa")

# --- 7. Move the _GoBack bookmark onto the empty Ghi chu cell of Stt=8 row (table row 9) ---
$targetCell = $t.Cell(9, 5)
$d.Bookmarks.Add("_GoBack", $targetCell.Range) | Out-Null
